$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 471
$ws.Range("I6").Value = 40
$ws.Range("J6").Value = 902
$ws.Range("K6").Value = 120
$ws.Range("L6").Value = 2706
$ws.Range("M6").Value = -8
$ws.Range("N6").Value = -2930
$ws.Range("H31").Value = 1098
$ws.Range("I31").Value = 1098
$ws.Range("K31").Value = 3294
$ws.Range("M31").Value = -3064
$ws.Range("H38").Value = 371.5
$ws.Range("I38").Value = 73
$ws.Range("J38").Value = 670
$ws.Range("K38").Value = 219
$ws.Range("L38").Value = 2010
$ws.Range("M38").Value = 153
$ws.Range("N38").Value = -2754
$ws.Range("H39").Value = 1063.375
$ws.Range("I39").Value = 152
$ws.Range("J39").Value = 1974.75
$ws.Range("K39").Value = 456
$ws.Range("L39").Value = 5924.25
$ws.Range("M39").Value = -160
$ws.Range("N39").Value = -6516.25
$ws.Range("H112").Value = 1385.1724
$ws.Range("J112").Value = 1423.8462
$ws.Range("L112").Value = 4271.5386
$ws.Range("N112").Value = -6487.5386
$ws.Range("H129").Value = 1472.1
$ws.Range("J129").Value = 1575
$ws.Range("L129").Value = 4725
$ws.Range("N129").Value = -14725
$ws.Range("H132").Value = 558730.4
$ws.Range("I132").Value = 4260.5454
$ws.Range("J132").Value = 1430040.1
$ws.Range("K132").Value = 12781.6362
$ws.Range("L132").Value = 4290120.300000001
$ws.Range("M132").Value = -10251.6362
$ws.Range("N132").Value = -4295180.300000001
$ws.Range("H137").Value = 1847.4286
$ws.Range("I137").Value = 1586.4
$ws.Range("K137").Value = 4759.200000000001
$ws.Range("M137").Value = -2209.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1737.9231
$ws.Range("I45").Value = 1284.7142
$ws.Range("J45").Value = 2266.6667
$ws.Range("K45").Value = 1284.7142
$ws.Range("L45").Value = 2266.6667
$ws.Range("M45").Value = -907.7141999999999
$ws.Range("N45").Value = -3020.6667
$ws.Range("H61").Value = 1923.3529
$ws.Range("I61").Value = 1978.4286
$ws.Range("J61").Value = 1666.3334
$ws.Range("K61").Value = 1978.4286
$ws.Range("L61").Value = 1666.3334
$ws.Range("M61").Value = -1766.4286
$ws.Range("N61").Value = -2090.3334
$ws.Range("H74").Value = 673.375
$ws.Range("I74").Value = 606.9
$ws.Range("J74").Value = 872.8
$ws.Range("K74").Value = 606.9
$ws.Range("L74").Value = 872.8
$ws.Range("M74").Value = 267.1
$ws.Range("N74").Value = -2620.8
$ws.Range("H77").Value = 673.375
$ws.Range("I77").Value = 606.9
$ws.Range("J77").Value = 872.8
$ws.Range("K77").Value = 3034.5
$ws.Range("L77").Value = 4364
$ws.Range("M77").Value = 1333.5
$ws.Range("N77").Value = -13100
$ws.Range("H109").Value = 20000
$ws.Range("J109").Value = 20000
$ws.Range("L109").Value = 20000
$ws.Range("N109").Value = -22774
$ws.Range("H136").Value = 1923.3529
$ws.Range("I136").Value = 1978.4286
$ws.Range("J136").Value = 1666.3334
$ws.Range("K136").Value = 5935.2858
$ws.Range("L136").Value = 4999.0002
$ws.Range("M136").Value = -3385.2858
$ws.Range("N136").Value = -10099.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 13000
$ws.Range("J38").Value = 13000
$ws.Range("L38").Value = 13000
$ws.Range("N38").Value = -13832
$ws.Range("H104").Value = 48500
$ws.Range("J104").Value = 48500
$ws.Range("L104").Value = 48500
$ws.Range("N104").Value = -55488
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 281.2
$ws.Range("I22").Value = 263.69232
$ws.Range("J22").Value = 395
$ws.Range("K22").Value = 263.69232
$ws.Range("L22").Value = 395
$ws.Range("M22").Value = 86.30768
$ws.Range("N22").Value = -1095
$ws.Range("H35").Value = 1148.5
$ws.Range("I35").Value = 864.6667
$ws.Range("J35").Value = 2000
$ws.Range("K35").Value = 864.6667
$ws.Range("L35").Value = 2000
$ws.Range("M35").Value = -570.6667
$ws.Range("N35").Value = -2588
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H94").Value = 4422.8184
$ws.Range("I94").Value = 4141.6665
$ws.Range("J94").Value = 4528.25
$ws.Range("K94").Value = 4141.6665
$ws.Range("L94").Value = 4528.25
$ws.Range("M94").Value = -3690.6665
$ws.Range("N94").Value = -5430.25
$ws.Range("H105").Value = 1240
$ws.Range("I105").Value = 1366.6666
$ws.Range("K105").Value = 1366.6666
$ws.Range("M105").Value = 380.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 7266.6665
$ws.Range("I98").Value = 900
$ws.Range("J98").Value = 10450
$ws.Range("K98").Value = 2700
$ws.Range("L98").Value = 31350
$ws.Range("M98").Value = -1202
$ws.Range("N98").Value = -34346
$ws.Range("H121").Value = 462
$ws.Range("I121").Value = 376.5
$ws.Range("J121").Value = 633
$ws.Range("K121").Value = 1129.5
$ws.Range("L121").Value = 1899
$ws.Range("M121").Value = 180.5
$ws.Range("N121").Value = -4519
$ws.Range("H129").Value = 1479.32
$ws.Range("I129").Value = 1290
$ws.Range("J129").Value = 1515.381
$ws.Range("K129").Value = 3870
$ws.Range("L129").Value = 4546.143
$ws.Range("M129").Value = 1130
$ws.Range("N129").Value = -14546.143

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 11666668
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H4").Value = 3000
$ws.Range("J4").Value = 3000
$ws.Range("L4").Value = 3000
$ws.Range("N4").Value = -3224
$ws.Range("H5").Value = 5202.778
$ws.Range("I5").Value = 293.33334
$ws.Range("J5").Value = 6184.6665
$ws.Range("K5").Value = 293.33334
$ws.Range("L5").Value = 6184.6665
$ws.Range("M5").Value = -181.33334
$ws.Range("N5").Value = -6408.6665
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H10").Value = 20000000
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H11").Value = 33666556
$ws.Range("I11").Value = 37750000
$ws.Range("J11").Value = 999000
$ws.Range("K11").Value = 37750000
$ws.Range("L11").Value = 999000
$ws.Range("M11").Value = -37749861
$ws.Range("N11").Value = -999278
$ws.Range("H13").Value = 368.33334
$ws.Range("I13").Value = 105
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 105
$ws.Range("L13").Value = 500
$ws.Range("M13").Value = 34
$ws.Range("N13").Value = -778

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 368.46667
$ws.Range("I16").Value = 368.46667
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 368.46667
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -198.46667
$ws.Range("N16").ClearContents()
$ws.Range("H68").Value = 2386.6667
$ws.Range("I68").Value = 2416.6667
$ws.Range("J68").Value = 2366.6667
$ws.Range("K68").Value = 2416.6667
$ws.Range("L68").Value = 2366.6667
$ws.Range("M68").Value = -1667.6667
$ws.Range("N68").Value = -3864.6667
$ws.Range("H71").Value = 2386.6667
$ws.Range("I71").Value = 2416.6667
$ws.Range("J71").Value = 2366.6667
$ws.Range("K71").Value = 12083.3335
$ws.Range("L71").Value = 11833.3335
$ws.Range("M71").Value = -8339.333500000001
$ws.Range("N71").Value = -19321.3335
$ws.Range("H82").Value = 2863.6365
$ws.Range("I82").Value = 2683.3333
$ws.Range("K82").Value = 2683.3333
$ws.Range("M82").Value = -2322.3333
$ws.Range("H85").Value = 2863.6365
$ws.Range("I85").Value = 2683.3333
$ws.Range("K85").Value = 2683.3333
$ws.Range("M85").Value = -1435.3333
$ws.Range("H122").Value = 5788.5356
$ws.Range("I122").Value = 7530.8237
$ws.Range("J122").Value = 3095.9092
$ws.Range("K122").Value = 22592.4711
$ws.Range("L122").Value = 9287.7276
$ws.Range("M122").Value = -20142.4711
$ws.Range("N122").Value = -14187.7276
$ws.Range("H132").Value = 2091.4285
$ws.Range("I132").Value = 1958.25
$ws.Range("J132").Value = 3156.8572
$ws.Range("K132").Value = 5874.75
$ws.Range("L132").Value = 9470.571599999999
$ws.Range("M132").Value = -3344.75
$ws.Range("N132").Value = -14530.5716
$ws.Range("H136").Value = 2551.1562
$ws.Range("I136").Value = 1405.4286
$ws.Range("K136").Value = 4216.2858
$ws.Range("M136").Value = -1666.2858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1497.75
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
